$d = $word.ActiveDocument

# The document header contains a "PAGE" field whose cached/displayed
# result text is stale ("1"). Update it to reflect the document's
# current page count (5), matching docProps/app.xml's <Pages>5</Pages>
# and the target revision of the header.
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    $flds = $hdr.Range.Fields
    for ($i = 1; $i -le $flds.Count; $i++) {
        $f = $flds.Item($i)
        if ($f.Type -eq 33) {
            $res = $f.Result
            if ($res.Text -eq "1") {
                $res.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "5", 2) | Out-Null
            }
        }
    }
}
